$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.047128063602184
$ws.Range("D2").Value = 1.053276440205262
$ws.Range("E2").Value = 1.058950772824788
$ws.Range("F2").Value = 1.066364642304763
$ws.Range("I2").Value = 1.045649976774784
$ws.Range("J2").Value = 1.052178781974173
$ws.Range("K2").Value = 1.056022855291684
$ws.Range("L2").Value = 1.06168160668493
$ws.Range("M2").Value = 1.069075388704568
$ws.Range("N2").Value = 1.021223028434477

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.048048835748804
$ws.Range("D3").Value = 1.053992202943841
$ws.Range("E3").Value = 1.059848817301546
$ws.Range("F3").Value = 1.067237896132549
$ws.Range("I3").Value = 1.045887372272664
$ws.Range("J3").Value = 1.052748042992124
$ws.Range("K3").Value = 1.056551879775694
$ws.Range("L3").Value = 1.062393572418461
$ws.Range("M3").Value = 1.069764076728424
$ws.Range("N3").Value = 1.021414912316752

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.048645100884352
$ws.Range("D4").Value = 1.054455803149088
$ws.Range("E4").Value = 1.060431165663696
$ws.Range("F4").Value = 1.067803897793608
$ws.Range("I4").Value = 1.046040088301862
$ws.Range("J4").Value = 1.053116207775297
$ws.Range("K4").Value = 1.056893973834551
$ws.Range("L4").Value = 1.062854876975676
$ws.Range("M4").Value = 1.070210007805282
$ws.Range("N4").Value = 1.021538946011491

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.048895880642617
$ws.Range("D5").Value = 1.05465080784866
$ws.Range("E5").Value = 1.060676283391242
$ws.Range("F5").Value = 1.068042070249368
$ws.Range("I5").Value = 1.046104075314892
$ws.Range("J5").Value = 1.053270938862513
$ws.Range("K5").Value = 1.057037736386535
$ws.Range("L5").Value = 1.063048955322093
$ws.Range("M5").Value = 1.070397548650296
$ws.Range("N5").Value = 1.021591058705289

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.048937994077517
$ws.Range("D6").Value = 1.054683556237462
$ws.Range("E6").Value = 1.060717457191908
$ws.Range("F6").Value = 1.068082073621126
$ws.Range("I6").Value = 1.046114806402094
$ws.Range("J6").Value = 1.053296916200591
$ws.Range("K6").Value = 1.05706187156582
$ws.Range("L6").Value = 1.063081550447485
$ws.Range("M6").Value = 1.070429041732953
$ws.Range("N6").Value = 1.021599806821456

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.048648451385412
$ws.Range("D7").Value = 1.054458408392165
$ws.Range("E7").Value = 1.060434439767993
$ws.Range("F7").Value = 1.067807079382129
$ws.Range("I7").Value = 1.046040944144252
$ws.Range("J7").Value = 1.053118275479322
$ws.Range("K7").Value = 1.056895895009521
$ws.Range("L7").Value = 1.062857469688401
$ws.Range("M7").Value = 1.070212513456259
$ws.Range("N7").Value = 1.021539642466479

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.047439146075475
$ws.Range("D8").Value = 1.053518240711272
$ws.Range("E8").Value = 1.059254010631443
$ws.Range("F8").Value = 1.066659565260305
$ws.Range("I8").Value = 1.045730390590567
$ws.Range("J8").Value = 1.052371204101516
$ws.Range("K8").Value = 1.056201686525585
$ws.Range("L8").Value = 1.061922090661009
$ws.Range("M8").Value = 1.069308070328678
$ws.Range("N8").Value = 1.02128790273958

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.045311801577685
$ws.Range("D9").Value = 1.051865086485061
$ws.Range("E9").Value = 1.057183616043261
$ws.Range("F9").Value = 1.064644830641512
$ws.Range("I9").Value = 1.045176329768417
$ws.Range("J9").Value = 1.051053397628937
$ws.Range("K9").Value = 1.05497676428197
$ws.Range("L9").Value = 1.06027860412098
$ws.Range("M9").Value = 1.067716712020457
$ws.Range("N9").Value = 1.020843343608661

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043896066068527
$ws.Range("D10").Value = 1.050765452523736
$ws.Range("E10").Value = 1.055809952666697
$ws.Range("F10").Value = 1.063306699174576
$ws.Range("I10").Value = 1.044802402533229
$ws.Range("J10").Value = 1.050174000746207
$ws.Range("K10").Value = 1.054159112309153
$ws.Range("L10").Value = 1.059186234325634
$ws.Range("M10").Value = 1.066657489567466
$ws.Range("N10").Value = 1.020546349582252

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.043283642937219
$ws.Range("D11").Value = 1.050289902889894
$ws.Range("E11").Value = 1.055216726896504
$ws.Range("F11").Value = 1.062728485342939
$ws.Range("I11").Value = 1.044639415406012
$ws.Range("J11").Value = 1.049793021459708
$ws.Range("K11").Value = 1.053804827418134
$ws.Range("L11").Value = 1.058714022834133
$ws.Range("M11").Value = 1.066199250275334
$ws.Range("N11").Value = 1.020417605805126

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.043056252933451
$ws.Range("D12").Value = 1.050113354069208
$ws.Range("E12").Value = 1.054996615372269
$ws.Range("F12").Value = 1.062513893762311
$ws.Range("I12").Value = 1.044578713927997
$ws.Range("J12").Value = 1.049651480575929
$ws.Range("K12").Value = 1.053673195811557
$ws.Range("L12").Value = 1.058538742774472
$ws.Range("M12").Value = 1.066029102866137
$ws.Range("N12").Value = 1.020369763505734

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043105024715441
$ws.Range("D13").Value = 1.050151220226885
$ws.Range("E13").Value = 1.055043819202627
$ws.Range("F13").Value = 1.062559916087775
$ws.Range("I13").Value = 1.044591741868383
$ws.Range("J13").Value = 1.049681842836174
$ws.Range("K13").Value = 1.053701432778677
$ws.Range("L13").Value = 1.05857633547616
$ws.Range("M13").Value = 1.066065597187926
$ws.Range("N13").Value = 1.020380026790766

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.043264844941314
$ws.Range("D14").Value = 1.050275307434461
$ws.Range("E14").Value = 1.055198527525492
$ws.Range("F14").Value = 1.062710743397969
$ws.Range("I14").Value = 1.044634401079854
$ws.Range("J14").Value = 1.049781322215699
$ws.Range("K14").Value = 1.053793947404556
$ws.Range("L14").Value = 1.058699531667558
$ws.Range("M14").Value = 1.066185184543247
$ws.Range("N14").Value = 1.020413651574334

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.043363327577702
$ws.Range("D15").Value = 1.0503517738255
$ws.Range("E15").Value = 1.055293880139983
$ws.Range("F15").Value = 1.062803697350053
$ws.Range("I15").Value = 1.044660663537532
$ws.Range("J15").Value = 1.049842611036856
$ws.Range("K15").Value = 1.053850944189261
$ws.Range("L15").Value = 1.05877545289014
$ws.Range("M15").Value = 1.066258874656125
$ws.Range("N15").Value = 1.020434366126012

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043936723304899
$ws.Range("D16").Value = 1.050797025935587
$ws.Range("E16").Value = 1.055849356525673
$ws.Range("F16").Value = 1.063345098866675
$ws.Range("I16").Value = 1.044813196870015
$ws.Range("J16").Value = 1.050199281073705
$ws.Range("K16").Value = 1.054182620144855
$ws.Range("L16").Value = 1.059217590245348
$ws.Range("M16").Value = 1.066687910172108
$ws.Range("N16").Value = 1.020554890896931

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.044296560637186
$ws.Range("D17").Value = 1.051076482265631
$ws.Range("E17").Value = 1.056198216225834
$ws.Range("F17").Value = 1.063685029727033
$ws.Range("I17").Value = 1.044908589822852
$ws.Range("J17").Value = 1.050422959380725
$ws.Range("K17").Value = 1.054390609222799
$ws.Range("L17").Value = 1.059495144276413
$ws.Range("M17").Value = 1.066957143786372
$ws.Range("N17").Value = 1.020630454802242

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.044506505450858
$ws.Range("D18").Value = 1.051239542123159
$ws.Range("E18").Value = 1.056401852479511
$ws.Range("F18").Value = 1.063883421893159
$ws.Range("I18").Value = 1.044964127237788
$ws.Range("J18").Value = 1.05055340830455
$ws.Range("K18").Value = 1.054511902872767
$ws.Range("L18").Value = 1.059657113157281
$ws.Range("M18").Value = 1.067114222767601
$ws.Range("N18").Value = 1.020674516085437

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.044578100974349
$ws.Range("D19").Value = 1.051295151064177
$ws.Range("E19").Value = 1.056471312961437
$ws.Range("F19").Value = 1.063951088163502
$ws.Range("I19").Value = 1.044983046457461
$ws.Range("J19").Value = 1.050597884788236
$ws.Range("K19").Value = 1.054553256963125
$ws.Range("L19").Value = 1.059712353262564
$ws.Range("M19").Value = 1.067167789312515
$ws.Range("N19").Value = 1.02068953748006

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.044257947511748
$ws.Range("D20").Value = 1.051046493266476
$ws.Range("E20").Value = 1.056160771117432
$ws.Range("F20").Value = 1.063648546327282
$ws.Range("I20").Value = 1.044898365781781
$ws.Range("J20").Value = 1.050398962757392
$ws.Range("K20").Value = 1.054368296320001
$ws.Range("L20").Value = 1.059465357452743
$ws.Range("M20").Value = 1.066928253465239
$ws.Range("N20").Value = 1.020622348941373

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.043217779341204
$ws.Range("D21").Value = 1.050238764307283
$ws.Range("E21").Value = 1.055152963171374
$ws.Range("F21").Value = 1.062666323455831
$ws.Range("I21").Value = 1.044621843439583
$ws.Range("J21").Value = 1.049752028788458
$ws.Range("K21").Value = 1.053766705094825
$ws.Range("L21").Value = 1.058663250127695
$ws.Range("M21").Value = 1.066149967296206
$ws.Range("N21").Value = 1.020403750493546

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.042564312630186
$ws.Range("D22").Value = 1.049731443098704
$ws.Range("E22").Value = 1.054520697559465
$ws.Range("F22").Value = 1.062049819264731
$ws.Range("I22").Value = 1.044447052806166
$ws.Range("J22").Value = 1.049345113000349
$ws.Range("K22").Value = 1.053388262177402
$ws.Range("L22").Value = 1.058159630182508
$ws.Range("M22").Value = 1.065660993430693
$ws.Range("N22").Value = 1.020266186938042

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.042910677245587
$ws.Range("D23").Value = 1.050000332935371
$ws.Range("E23").Value = 1.054855741862123
$ws.Range("F23").Value = 1.062376538948953
$ws.Range("I23").Value = 1.044539800592984
$ws.Range("J23").Value = 1.049560841839384
$ws.Range("K23").Value = 1.053588900401874
$ws.Range("L23").Value = 1.058426542148145
$ws.Range("M23").Value = 1.065920172639709
$ws.Range("N23").Value = 1.020339123409991

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.044275394942597
$ws.Range("D24").Value = 1.051060043826194
$ws.Range("E24").Value = 1.056177690479745
$ws.Range("F24").Value = 1.063665031246222
$ws.Range("I24").Value = 1.044902985906265
$ws.Range("J24").Value = 1.050409805856896
$ws.Range("K24").Value = 1.054378378631075
$ws.Range("L24").Value = 1.059478816600568
$ws.Range("M24").Value = 1.066941307635673
$ws.Range("N24").Value = 1.020626011673871

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.045861336633009
$ws.Range("D25").Value = 1.052292037758077
$ws.Range("E25").Value = 1.057717706505466
$ws.Range("F25").Value = 1.0651648091432
$ws.Range("I25").Value = 1.045320372970777
$ws.Range("J25").Value = 1.051394238235678
$ws.Range("K25").Value = 1.055293623559916
$ws.Range("L25").Value = 1.05857633547616
$ws.Range("M25").Value = 1.066065597187926
$ws.Range("N25").Value = 1.020380026790766
